# Update "想去人数" (interest count) figures in column F across the
# workbook's sheets, reflecting newly generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 1303
$ws.Range("F7").Value  = 1011
$ws.Range("F8").Value  = 952
$ws.Range("F12").Value = 458
$ws.Range("F15").Value = 4361
$ws.Range("F16").Value = 1266
$ws.Range("F18").Value = 2755
$ws.Range("F19").Value = 679
$ws.Range("F22").Value = 3768
$ws.Range("F23").Value = 821
$ws.Range("F27").Value = 2482
$ws.Range("F29").Value = 891
$ws.Range("F31").Value = 987
$ws.Range("F32").Value = 254
$ws.Range("F33").Value = 9
$ws.Range("F35").Value = 56
$ws.Range("F36").Value = 1441
$ws.Range("F37").Value = 2014
$ws.Range("F38").Value = 956
$ws.Range("F39").Value = 11
$ws.Range("F40").Value = 5
$ws.Range("F41").Value = 528
$ws.Range("F42").Value = 107
$ws.Range("F45").Value = 309
$ws.Range("F47").Value = 176

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value  = 155

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value  = 498

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 498
$ws.Range("F4").Value  = 1303
$ws.Range("F6").Value  = 1011
$ws.Range("F7").Value  = 952
$ws.Range("F8").Value  = 155
$ws.Range("F15").Value = 458
$ws.Range("F17").Value = 4361
$ws.Range("F18").Value = 1266
$ws.Range("F21").Value = 2755
$ws.Range("F23").Value = 3768
$ws.Range("F24").Value = 821
$ws.Range("F28").Value = 2482
$ws.Range("F34").Value = 891
$ws.Range("F36").Value = 987
$ws.Range("F37").Value = 254
$ws.Range("F39").Value = 1441
$ws.Range("F40").Value = 2014
$ws.Range("F41").Value = 956
$ws.Range("F42").Value = 528
$ws.Range("F43").Value = 107
$ws.Range("F45").Value = 309
$ws.Range("F47").Value = 176
